$d = $word.ActiveDocument

$d.Content.Find.Execute("48×12=576", $true, $false, $false, $false, $false, $true, 1, $false, "61×30=1830", 2) | Out-Null
$d.Content.Find.Execute("92×35=3220", $true, $false, $false, $false, $false, $true, 1, $false, "98×60=5880", 2) | Out-Null
$d.Content.Find.Execute("79×87=6873", $true, $false, $false, $false, $false, $true, 1, $false, "32×63=2016", 2) | Out-Null
$d.Content.Find.Execute("74×71=5254", $true, $false, $false, $false, $false, $true, 1, $false, "85×96=8160", 2) | Out-Null
$d.Content.Find.Execute("95×89=8455", $true, $false, $false, $false, $false, $true, 1, $false, "88×23=2024", 2) | Out-Null
$d.Content.Find.Execute("79×41=3239", $true, $false, $false, $false, $false, $true, 1, $false, "15×33=495", 2) | Out-Null
$d.Content.Find.Execute("60×88=5280", $true, $false, $false, $false, $false, $true, 1, $false, "84×74=6216", 2) | Out-Null
$d.Content.Find.Execute("55×95=5225", $true, $false, $false, $false, $false, $true, 1, $false, "20×19=380", 2) | Out-Null
$d.Content.Find.Execute("25×71=1775", $true, $false, $false, $false, $false, $true, 1, $false, "18×17=306", 2) | Out-Null
$d.Content.Find.Execute("43×44=1892", $true, $false, $false, $false, $false, $true, 1, $false, "76×19=1444", 2) | Out-Null
$d.Content.Find.Execute("28×38=1064", $true, $false, $false, $false, $false, $true, 1, $false, "43×66=2838", 2) | Out-Null
$d.Content.Find.Execute("74×36=2664", $true, $false, $false, $false, $false, $true, 1, $false, "23×25=575", 2) | Out-Null
$d.Content.Find.Execute("99×84=8316", $true, $false, $false, $false, $false, $true, 1, $false, "64×49=3136", 2) | Out-Null
$d.Content.Find.Execute("87×24=2088", $true, $false, $false, $false, $false, $true, 1, $false, "27×93=2511", 2) | Out-Null
$d.Content.Find.Execute("58×82=4756", $true, $false, $false, $false, $false, $true, 1, $false, "59×41=2419", 2) | Out-Null
$d.Content.Find.Execute("32×33=1056", $true, $false, $false, $false, $false, $true, 1, $false, "85×87=7395", 2) | Out-Null
$d.Content.Find.Execute("37×20=740", $true, $false, $false, $false, $false, $true, 1, $false, "40×47=1880", 2) | Out-Null
$d.Content.Find.Execute("89×14=1246", $true, $false, $false, $false, $false, $true, 1, $false, "49×35=1715", 2) | Out-Null
$d.Content.Find.Execute("20×78=1560", $true, $false, $false, $false, $false, $true, 1, $false, "31×27=837", 2) | Out-Null
$d.Content.Find.Execute("60×52=3120", $true, $false, $false, $false, $false, $true, 1, $false, "42×78=3276", 2) | Out-Null
$d.Content.Find.Execute("16×14=224", $true, $false, $false, $false, $false, $true, 1, $false, "65×37=2405", 2) | Out-Null
$d.Content.Find.Execute("15×38=570", $true, $false, $false, $false, $false, $true, 1, $false, "75×79=5925", 2) | Out-Null
$d.Content.Find.Execute("52×66=3432", $true, $false, $false, $false, $false, $true, 1, $false, "81×73=5913", 2) | Out-Null
$d.Content.Find.Execute("87×77=6699", $true, $false, $false, $false, $false, $true, 1, $false, "75×39=2925", 2) | Out-Null
$d.Content.Find.Execute("51×55=2805", $true, $false, $false, $false, $false, $true, 1, $false, "39×24=936", 2) | Out-Null
